# Updates after call with GetWireless (2016-12-14)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issue Tracking")

# Update Status column (G) for issue rows 11-13
$ws.Range("G11").Value = "Under Test (GW)"
$ws.Range("G12").Value = "Under Test (OP)"
$ws.Range("G13").Value = "Under Test (OP)"

# Update the frozen-pane view state (scroll position) and active selection
$ws.Activate()
$ws.Range("A3").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F25").Select()
